$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace spaces with underscores in the "experimentDesign" (col D) and
# "strain" (col F) data columns so the values conform to spec.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    if ($dCell.Value2 -eq "Environmental Perturbation") {
        $dCell.Value = "Environmental_Perturbation"
    }

    $fCell = $ws.Cells.Item($r, 6)
    if ($fCell.Value2 -eq "KN99 alpha") {
        $fCell.Value = "KN99_alpha"
    }
}

# Update the active selection on the sheet from B2:B27 to F2:F27
$ws.Range("F2:F27").Select()
